# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns,
# and correct the B/C/D/E ranking-row data for two swapped coin pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.245.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.439.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.28%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.27%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.442.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.12%  "

$ws.Range("E10").Value = "  -5.57%  "

$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.884.40"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.182.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.440.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "317.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.95%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0975"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.568.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.77%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "533.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.70%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.81%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.67%  "

$ws.Range("E33").Value = "  -5.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.72%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.377"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.93%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.24%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "141.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0532"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.587"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0928"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.54%  "
